$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the preparer name "Retrofitted_1476" with "S.GISH" in the
# s2cDNAPreparer (B) and libraryPreparer (E) columns for data rows 2-7.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 2).Value = "S.GISH"
    $ws.Cells.Item($row, 5).Value = "S.GISH"
}

# Update the active selection on the sheet from N9 to B9.
$ws.Range("B9").Select()
